$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    Touches the Overview sheet (columns E & F, rows 2-3) and the locale
#    sheets' "Status" column (C, rows 2-3).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. Fill in the handback details for the zh-cn sheet: the "Latest Target
#    File" (I) becomes a hyperlink to the source .md file, "Latest Handback
#    File" (J) gets the generated .xlf name, and "Latest Handback DateTime"
#    (K) records when the handback xliff was produced.
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d56fb39dd30b3b67d462bb71a0f0b46d4d12bd6a/e2e/b86e030c-1c66-4fef-abf0-d3b1b3c5cbdb.md", "", "", "b86e030c-1c66-4fef-abf0-d3b1b3c5cbdb.md") | Out-Null
$wsZhCn.Range("J2").Value = "b86e030c-1c66-4fef-abf0-d3b1b3c5cbdb.3e9019e7ebf5c452300091e5259f1578bae13eab.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-24 18:31:50"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d56fb39dd30b3b67d462bb71a0f0b46d4d12bd6a/e2e/cfdd0013-44eb-4700-b218-01af24f046ee.md", "", "", "cfdd0013-44eb-4700-b218-01af24f046ee.md") | Out-Null
$wsZhCn.Range("J3").Value = "cfdd0013-44eb-4700-b218-01af24f046ee.135a840f96fc75608b792cd829fbf2a95beddc0f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-24 18:31:50"

# ---------------------------------------------------------------------------
# 3. Same handback bookkeeping for the de-de sheet.
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d56fb39dd30b3b67d462bb71a0f0b46d4d12bd6a/e2e/b86e030c-1c66-4fef-abf0-d3b1b3c5cbdb.md", "", "", "b86e030c-1c66-4fef-abf0-d3b1b3c5cbdb.md") | Out-Null
$wsDeDe.Range("J2").Value = "b86e030c-1c66-4fef-abf0-d3b1b3c5cbdb.3e9019e7ebf5c452300091e5259f1578bae13eab.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-24 18:31:57"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d56fb39dd30b3b67d462bb71a0f0b46d4d12bd6a/e2e/cfdd0013-44eb-4700-b218-01af24f046ee.md", "", "", "cfdd0013-44eb-4700-b218-01af24f046ee.md") | Out-Null
$wsDeDe.Range("J3").Value = "cfdd0013-44eb-4700-b218-01af24f046ee.135a840f96fc75608b792cd829fbf2a95beddc0f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-24 18:31:57"

# ---------------------------------------------------------------------------
# 4. Widen the columns that now carry longer hyperlink / file-name text.
# ---------------------------------------------------------------------------
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 29.9777047293527
$wsZhCn.Range("I1").EntireColumn.ColumnWidth = 40
$wsZhCn.Range("J1").EntireColumn.ColumnWidth = 40

$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 29.9777047293527
$wsDeDe.Range("I1").EntireColumn.ColumnWidth = 40
$wsDeDe.Range("J1").EntireColumn.ColumnWidth = 40

Write-Host "Handback report generated"
